$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures as literal text (e.g. "235.01", "29.259.30").
# Force a Text number format before writing so Excel does not silently
# reinterpret these strings as numbers and drop trailing zeros / precision.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.259.30"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.88"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.01"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6014"
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06991"
$ws.Range("E8").Value = "  -5.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2772"
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.45"
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07641"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.829.19"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6262"
$ws.Range("E14").Value = "  -7.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009721"
$ws.Range("E15").Value = "  -5.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.71"
$ws.Range("E16").Value = "  -3.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.203.31"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.763"
$ws.Range("E18").Value = "  -7.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "222.54"
$ws.Range("E19").Value = "  -5.16%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.60"
$ws.Range("E21").Value = "  -5.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.920"
$ws.Range("E22").Value = "  -5.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.35"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1302"
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.970"
$ws.Range("E26").Value = "  -6.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.56"
$ws.Range("E27").Value = "  -4.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06944"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.454"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.444"
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.836"
$ws.Range("E31").Value = "  -4.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.773"
$ws.Range("E32").Value = "  -7.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.099"
$ws.Range("E33").Value = "  -3.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.724"
$ws.Range("E34").Value = "  -5.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6450"
$ws.Range("E35").Value = "  -8.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.543"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.742"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.203.86"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01743"
$ws.Range("E39").Value = "  -5.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.508"
$ws.Range("E40").Value = "  -5.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9041"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.986.25"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.30"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.26"
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("E46").Value = "  -4.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.488"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.578"
$ws.Range("E48").Value = "  -7.38%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4553"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05507"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.407"
$ws.Range("E51").Value = "  -7.67%  "
